$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.136775493621826
$ws.Range("B1").Value = 2.50940990447998
$ws.Range("C1").Value = 3.990776062011719
$ws.Range("D1").Value = 3.586139678955078
$ws.Range("E1").Value = 1.232507705688477
